# Weekly update: insert two new price records at the top of the data table
# (row 22), pushing all existing records below down by two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 22 (each Insert() call pushes rows 22+ down by one)
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

# New record data (identical for both inserted rows, per source data)
$mercado    = "Agrícola del Norte S.A. de Arica"
$region     = "Arica y Parinacota"
$fecha      = 44923
$codreg     = 15
$catId      = 100112045
$categoria  = "Zapallo"
$variedad   = "Camote"
$calidad    = "1a nueva(o)"
$volumen    = 500
$precioMin  = 850
$precioMax  = 900
$precioProm = 865
$unidad     = "`$/kilo (volumen en unidades)"
$origen     = "Perú"
$precioKg   = 865
$kgOUnidad  = 1
$clasif     = "Hortaliza"

# Column D (Fecha) keeps the same custom date number-format used by the rest of the table
$fechaFormat = $ws.Cells.Item(25, 4).NumberFormat

foreach ($r in 22, 23) {
    $ws.Cells.Item($r, 1).Value = 1
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).NumberFormat = $fechaFormat
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $catId
    $ws.Cells.Item($r, 7).Value = $categoria
    $ws.Cells.Item($r, 8).Value = $variedad
    $ws.Cells.Item($r, 9).Value = $calidad
    $ws.Cells.Item($r, 10).Value = $volumen
    $ws.Cells.Item($r, 11).Value = $precioMin
    $ws.Cells.Item($r, 12).Value = $precioMax
    $ws.Cells.Item($r, 13).Value = $precioProm
    $ws.Cells.Item($r, 14).Value = $unidad
    $ws.Cells.Item($r, 15).Value = $origen
    $ws.Cells.Item($r, 16).Value = $precioKg
    $ws.Cells.Item($r, 17).Value = $kgOUnidad
    $ws.Cells.Item($r, 18).Value = $clasif
}
